# Rename the worksheet from "PRECEPTORS" to "Preceptor Schedule"
# and update the active selection from C115 to F19 (cosmetic view-state
# changes captured by the commit: switching the file's display name and
# leaving the cursor on the cell that was last worked on).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Name = "Preceptor Schedule"

$ws.Activate()
$ws.Range("F19").Select()
